$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price text (e.g. "1.001", "29.947.43") that must
# stay literal text (it used "." as a thousands separator), so force those specific
# cells to Text format before assigning - otherwise Excel silently reinterprets
# strings like "1.001" or "5.400" as numbers and drops the trailing zero.
$ws.Range("D2:D22").NumberFormat = "@"
$ws.Range("D25:D28").NumberFormat = "@"
$ws.Range("D31:D36").NumberFormat = "@"
$ws.Range("D39:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.947.43"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "1.892.61"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "0.7764"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").Value = "243.92"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.3130"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").Value = "25.80"
$ws.Range("E9").Value = "  +2.31%  "

$ws.Range("D10").Value = "0.07245"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("D11").Value = "0.08704"

$ws.Range("D12").Value = "2.087.04"
$ws.Range("E12").Value = "  +7.83%  "

$ws.Range("D13").Value = "0.7727"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").Value = "5.400"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").Value = "94.39"
$ws.Range("E15").Value = "  +2.35%  "

$ws.Range("D16").Value = "6.197"
$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("D17").Value = "30.198.19"
$ws.Range("E17").Value = "  +1.35%  "

$ws.Range("D18").Value = "13.91"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").Value = "2.319.55"
$ws.Range("E19").Value = "  +9.44%  "

$ws.Range("D20").Value = "245.85"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").Value = "0.000007865"
$ws.Range("E21").Value = "  +1.57%  "

$ws.Range("D22").Value = "8.151"
$ws.Range("E22").Value = "  +0.87%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "0.1647"
$ws.Range("E25").Value = "  +4.24%  "

$ws.Range("D26").Value = "9.495"
$ws.Range("E26").Value = "  +1.29%  "

$ws.Range("D27").Value = "162.93"
$ws.Range("E27").Value = "  +0.71%  "

$ws.Range("D28").Value = "18.83"
$ws.Range("E28").Value = "  +0.60%  "

$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").Value = "1.542"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D32").Value = "4.514"
$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("D33").Value = "4.123"
$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").Value = "0.05478"
$ws.Range("E34").Value = "  -0.73%  "

$ws.Range("D35").Value = "1.245"
$ws.Range("E35").Value = "  -1.07%  "

$ws.Range("D36").Value = "0.7522"
$ws.Range("E36").Value = "  +1.02%  "

$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").Value = "0.01969"
$ws.Range("E39").Value = "  +3.01%  "

$ws.Range("D40").Value = "2.788"
$ws.Range("E40").Value = "  +0.36%  "

$ws.Range("D41").Value = "0.4501"
$ws.Range("E41").Value = "  +2.06%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.106.41"
$ws.Range("E42").Value = "  -2.84%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "73.65"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("D44").Value = "6.094"
$ws.Range("E44").Value = "  +3.90%  "

$ws.Range("D45").Value = "0.8539"
$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("D46").Value = "2.200.88"
$ws.Range("E46").Value = "  +8.60%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").Value = "103.29"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").Value = "1.880"
$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.909"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "7.600"
$ws.Range("E51").Value = "  +2.27%  "
